$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the element rows (2-7) so the elements appear in the order
# Ti, V, Nb, Hf, Si, Cr -- matching the order elements first appear in
# the alloy composition column headers (row 1), keeping each element's
# measured values together with its label.

$data = @(
    @("Ti", 2.09, 1.32, 1.21, 2.09, 1.1),
    @("V",  0.83, 1.4,  1.28, 0.87, 1.21),
    @("Nb", 2.28, 1.46, 2.29, 2.5,  2.09),
    @("Hf", 4.74, 4.8,  4.52, 4.67, 4.13),
    @("Si", 0,    0,    1.15, 0.92, 0),
    @("Cr", 0,    0,    0,    0,    1.91)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
}

$ws.Range("C14").Select()

# The sheet's table was originally mis-numbered ("Table2") from the
# openpyxl-authored workbook; Excel normalises this to "Table1" once the
# file is touched/re-saved through the UI.
$tbl = $ws.ListObjects.Item(1)
$tbl.Name = "Table1"
